# Bump the published ValueSet metadata to version 1.1.0, matching the
# commit "Added 1.1.0 of term": update the Version and Date rows on the
# Metadata sheet's Property/Value table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$versionRow = $ws.Columns("A").Find("Version")
if ($versionRow -ne $null) {
    $ws.Cells.Item($versionRow.Row, 2).Value = "1.1.0"
} else {
    $ws.Range("B3").Value = "1.1.0"
}

$dateRow = $ws.Columns("A").Find("Date")
if ($dateRow -ne $null) {
    $ws.Cells.Item($dateRow.Row, 2).Value = "2023-07-10T23:08:03+02:00"
} else {
    $ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
}
